$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Development row for 2021-12-15 (serial 44545)
$ws.Range("A3").Value = 44545
$ws.Range("A3").NumberFormat = "d-mmm-yy"
$ws.Range("B3").Value = "Development"
$ws.Range("C3").Value = 119
$ws.Range("D3").Value = 119
$ws.Range("E3").Value = 0

# Row 4: Production row for 2021-12-15 (serial 44545)
$ws.Range("A4").Value = 44545
$ws.Range("A4").NumberFormat = "d-mmm-yy"
$ws.Range("B4").Value = "Production"
$ws.Range("C4").Value = 134
$ws.Range("D4").Value = 132
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = "After execution all test cases pass"
$ws.Range("G4").Value = "Test cases iitially fail because of page load affected by network"

$ws.Rows.Item(4).RowHeight = 75

$ws.Range("F4:G4").Select()
